# The "setHeightConserveRatio" template had its image-expression field
# (" m:'Mona_Lisa.jpg'.asImage().setHeight(100) ", stored as a real Word
# field with fldChar begin/end + w:instrText runs) rewritten as plain
# literal text runs using { and } to delimit the M2Doc expression, per
# the move to TokenIteratorFieldRewriterSplit. Convert the field in place:
# drop the field's begin/end fldChars and its leading/trailing instrText
# space runs, and turn every remaining w:instrText run into an equivalent
# w:t run (preserving each run's rPr, e.g. the orange color), keeping the
# _GoBack bookmark where it sits among the runs.

$d = $word.ActiveDocument

$field = $d.Fields(1)

$rangeStart = $field.Code.Start - 1
$rangeEnd = $field.Result.End

$target = $d.Range($rangeStart, $rangeEnd)

$color = '<w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr>'

$runs = ""
$runs += "<w:r><w:t>{</w:t></w:r>"
$runs += "<w:r><w:t>m</w:t></w:r>"
$runs += "<w:r><w:t>:</w:t></w:r>"
$runs += "<w:r>$color<w:t>'</w:t></w:r>"
$runs += "<w:r>$color<w:t>Mona_Lisa</w:t></w:r>"
$runs += "<w:r>$color<w:t>.jpg</w:t></w:r>"
$runs += "<w:r>$color<w:t>'.asImage()</w:t></w:r>"
$runs += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$runs += "<w:r>$color<w:t>.setHeight(100)</w:t></w:r>"
$runs += '<w:r><w:t xml:space="preserve">}</w:t></w:r>'

$xml = '<?xml version="1.0" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' + $runs + '</w:p></w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
